# petty-cashBook-2021.xlsx — 30-Jun-2021 midday update
# Adds new Buku KAS (Sheet1) petty-cash transactions for 28/29/30-Jun-2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- 28-Jun-2021 (row 3 already has date + "Wages Expense"; add its Debit) ---
$ws.Range("D3").Formula = "=60000+260000"

$ws.Range("B4").Value = "TRANSFER BCA"
$ws.Range("D4").Formula = "=1136000+3312000+190000+51960000+331000"

$ws.Range("B5").Value = "A/R"
$ws.Range("C5").Formula = "=44100000+7860000+9570500"

$ws.Range("B6").Value = "SALES - cash/retail"
$ws.Range("C6").Formula = "=4549775+9376725-9570500"

$ws.Range("B7").Value = "SELISIH - lebih"
$ws.Range("C7").Value = 440000

$ws.Range("B8").Value = "SETOR KE BANK"
$ws.Range("D8").Value = 9000000

# --- 29-Jun-2021 (row 9) ---
$ws.Range("A9").Value = 44376
$ws.Range("B9").Value = "Wages Expense"
$ws.Range("D9").Formula = "=60000+260000"

$ws.Range("B10").Value = "A/R"
$ws.Range("C10").Formula = "=91200000+6260000+30808000"

$ws.Range("B11").Value = "TRANSFER BCA"
$ws.Range("D11").Formula = "=91200000+1426000+312000+6260000+27610000"

$ws.Range("B12").Value = "SALES - cash/retail"
$ws.Range("C12").Formula = "=28851275+5426725-30808000"

$ws.Range("B13").Value = "SELISIH - lebih"
$ws.Range("C13").Value = 70000

$ws.Range("B14").Value = "SETOR KE BANK"
$ws.Range("D14").Formula = "=5000000"

# --- 30-Jun-2021 (row 15) ---
$ws.Range("A15").Value = 44377
$ws.Range("B15").Value = "Wages Expense"
$ws.Range("D15").Formula = "=60000"

$ws.Range("B16").Value = "TRANSFER BCA"
$ws.Range("D16").Formula = "=1200000+922000+6750000+246000+33000+3245000"

$ws.Range("B17").Value = "CHEQUE RECEIVED"
$ws.Range("D17").Formula = "=1800000"

$ws.Range("B18").Value = "A/P"
$ws.Range("D18").Formula = "=700000"

$ws.Range("B19").Value = "BELI kresek"
$ws.Range("D19").Value = 54000

$ws.Range("B20").Value = "BELI isi stapler"
$ws.Range("D20").Value = 25000

$ws.Range("B21").Value = "A/R"
$ws.Range("C21").Formula = "=6750000"

$ws.Range("B22").Value = "FREIGHT OUT"
$ws.Range("D22").Formula = "=43500"

# --- View: midday update landed on row 14, so the selection moved there ---
$ws.Range("A14").Select() | Out-Null
